# Excel Push 30th Aug
# - addListItem!A2: "ListItemQ" -> "ListItemR" (C2 formula =A2 recalculates)
# - addListItem!D2: "ADLILC.8850" -> "ADLILC.8851"
# - createUser!A2: 20 -> 21 (B2/F2 formulas recalculate)
# - active sheet moves from addListItem (tab 2) to createUser (tab 3), with
#   the selection on createUser landing on A2

$wb = $excel.ActiveWorkbook

$wsAddListItem = $wb.Worksheets.Item("addListItem")
$wsCreateUser  = $wb.Worksheets.Item("createUser")

$wsAddListItem.Range("A2").Value = "ListItemR"
$wsAddListItem.Range("D2").Value = "ADLILC.8851"

$wsCreateUser.Range("A2").Value = 21

$wsCreateUser.Activate()
$wsCreateUser.Range("A2").Select()
